# Update "想去人数" (want-to-go count) values across sheets to reflect
# the latest scrape output (gh-pages generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# -------- Sheet: 展览 (Exhibitions) --------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 1273
$ws1.Range("F6").Value  = 28
$ws1.Range("F10").Value = 3515
$ws1.Range("F11").Value = 136
$ws1.Range("F16").Value = 608
$ws1.Range("F17").Value = 98
$ws1.Range("F18").Value = 761
$ws1.Range("F20").Value = 127
$ws1.Range("F22").Value = 61
$ws1.Range("F24").Value = 2697
$ws1.Range("F25").Value = 5195
$ws1.Range("F27").Value = 76
$ws1.Range("F29").Value = 3082
$ws1.Range("F30").Value = 288
$ws1.Range("F31").Value = 2258
$ws1.Range("F34").Value = 84
$ws1.Range("F38").Value = 31
$ws1.Range("F43").Value = 453
$ws1.Range("F45").Value = 487

# -------- Sheet: 演出 (Performances) --------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 74

# -------- Sheet: 全部类型 (All types, combined) --------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 1273
$ws4.Range("F6").Value  = 28
$ws4.Range("F10").Value = 3515
$ws4.Range("F11").Value = 136
$ws4.Range("F14").Value = 74
$ws4.Range("F17").Value = 608
$ws4.Range("F18").Value = 98
$ws4.Range("F19").Value = 761
$ws4.Range("F21").Value = 127
$ws4.Range("F23").Value = 61
$ws4.Range("F25").Value = 2697
$ws4.Range("F26").Value = 5195
$ws4.Range("F28").Value = 76
$ws4.Range("F30").Value = 3082
$ws4.Range("F31").Value = 288
$ws4.Range("F32").Value = 2258
$ws4.Range("F35").Value = 84
$ws4.Range("F39").Value = 31
$ws4.Range("F44").Value = 453
$ws4.Range("F46").Value = 487
